$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell H1 "Save" - copy formatting from the neighboring header cell (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "Save" column values for rows 2-8
$values = @(0, 1, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
